# Auto-generated Excel COM-interop script
# Updates Leve profit calculation columns (H-N) across several sheets
# to reflect refreshed market-board price data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 241.28572
$ws.Range("I11").Value = 241.28572
$ws.Range("K11").Value = 241.28572
$ws.Range("M11").Value = -101.28572

# Row 31
$ws.Range("H31").Value = 4757.1
$ws.Range("I31").Value = 367.2857
$ws.Range("K31").Value = 1101.8571
$ws.Range("M31").Value = -871.8571000000002

# Row 39
$ws.Range("H39").Value = 153.16667
$ws.Range("I39").Value = 136
$ws.Range("J39").Value = 187.5
$ws.Range("K39").Value = 408
$ws.Range("L39").Value = 562.5
$ws.Range("M39").Value = -112
$ws.Range("N39").Value = -1154.5

# Row 80
$ws.Range("H80").Value = 4977.8887
$ws.Range("I80").Value = 2599.625
$ws.Range("J80").Value = 6880.5
$ws.Range("K80").Value = 7798.875
$ws.Range("L80").Value = 20641.5
$ws.Range("M80").Value = -6800.875
$ws.Range("N80").Value = -22637.5

# Row 83
$ws.Range("H83").Value = 4977.8887
$ws.Range("I83").Value = 2599.625
$ws.Range("J83").Value = 6880.5
$ws.Range("K83").Value = 23396.625
$ws.Range("L83").Value = 61924.5
$ws.Range("M83").Value = -18404.625
$ws.Range("N83").Value = -71908.5

# Row 88
$ws.Range("H88").Value = 1694.6471
$ws.Range("J88").Value = 1817.1666
$ws.Range("L88").Value = 1817.1666
$ws.Range("N88").Value = -2629.1666

# Row 91
$ws.Range("H91").Value = 1694.6471
$ws.Range("J91").Value = 1817.1666
$ws.Range("L91").Value = 1817.1666
$ws.Range("N91").Value = -4625.1666

# Row 98
$ws.Range("H98").Value = 2839.0625
$ws.Range("I98").Value = 1934.5
$ws.Range("J98").Value = 4346.6665
$ws.Range("K98").Value = 1934.5
$ws.Range("L98").Value = 4346.6665
$ws.Range("M98").Value = -436.5
$ws.Range("N98").Value = -7342.6665

# Row 122
$ws.Range("H122").Value = 2839.0625
$ws.Range("I122").Value = 1934.5
$ws.Range("J122").Value = 4346.6665
$ws.Range("K122").Value = 5803.5
$ws.Range("L122").Value = 13039.9995
$ws.Range("M122").Value = -3353.5
$ws.Range("N122").Value = -17939.9995

# Row 132
$ws.Range("H132").Value = 2382.75
$ws.Range("I132").Value = 1074.4348
$ws.Range("J132").Value = 4152.8237
$ws.Range("K132").Value = 3223.3044
$ws.Range("L132").Value = 12458.4711
$ws.Range("M132").Value = -693.3044
$ws.Range("N132").Value = -17518.4711

# Row 138
$ws.Range("H138").Value = 3511356.5
$ws.Range("I138").Value = 1490.0605
$ws.Range("J138").Value = 8337423
$ws.Range("K138").Value = 4470.181500000001
$ws.Range("L138").Value = 25012269
$ws.Range("M138").Value = 669.8184999999994
$ws.Range("N138").Value = -25022549

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 1875.75
$ws.Range("I88").Value = 1932
$ws.Range("J88").Value = 1707
$ws.Range("K88").Value = 1932
$ws.Range("L88").Value = 1707
$ws.Range("M88").Value = -1526
$ws.Range("N88").Value = -2519

# Row 91
$ws.Range("H91").Value = 1875.75
$ws.Range("I91").Value = 1932
$ws.Range("J91").Value = 1707
$ws.Range("K91").Value = 1932
$ws.Range("L91").Value = 1707
$ws.Range("M91").Value = -528
$ws.Range("N91").Value = -4515

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1750
$ws.Range("I5").Value = 1750
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1750
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1637
$ws.Range("N5").ClearContents()

# Row 99
$ws.Range("H99").Value = 2512.5386
$ws.Range("I99").Value = 2739.111
$ws.Range("J99").Value = 2002.75
$ws.Range("K99").Value = 2739.111
$ws.Range("L99").Value = 2002.75
$ws.Range("M99").Value = -1241.111
$ws.Range("N99").Value = -4998.75

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 74.916664
$ws.Range("I10").Value = 49.88889
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 149.66667
$ws.Range("L10").Value = 450
$ws.Range("M10").Value = -10.66667000000001
$ws.Range("N10").Value = -728

# Row 21
$ws.Range("H21").Value = 1754.1277
$ws.Range("I21").Value = 2582.8333
$ws.Range("J21").Value = 1470
$ws.Range("K21").Value = 7748.499899999999
$ws.Range("L21").Value = 4410
$ws.Range("M21").Value = -7575.499899999999
$ws.Range("N21").Value = -4756

# Row 23
$ws.Range("H23").Value = 445.55554
$ws.Range("I23").Value = 405
$ws.Range("J23").Value = 457.14285
$ws.Range("K23").Value = 1215
$ws.Range("L23").Value = 1371.42855
$ws.Range("M23").Value = -980
$ws.Range("N23").Value = -1841.42855

# Row 68
$ws.Range("H68").Value = 1354.9512
$ws.Range("I68").Value = 685.3333
$ws.Range("J68").Value = 2300.2942
$ws.Range("K68").Value = 2055.9999
$ws.Range("L68").Value = 6900.882599999999
$ws.Range("M68").Value = -1244.9999
$ws.Range("N68").Value = -8522.882599999999

# Row 71
$ws.Range("H71").Value = 1354.9512
$ws.Range("I71").Value = 685.3333
$ws.Range("J71").Value = 2300.2942
$ws.Range("K71").Value = 6167.9997
$ws.Range("L71").Value = 20702.6478
$ws.Range("M71").Value = -2111.9997
$ws.Range("N71").Value = -28814.6478

# Row 107
$ws.Range("H107").Value = 504586.75
$ws.Range("I107").Value = 516.86957
$ws.Range("K107").Value = 1550.60871
$ws.Range("M107").Value = 369.39129

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3353.75
$ws.Range("I7").Value = 2563.3333
$ws.Range("K7").Value = 2563.3333
$ws.Range("M7").Value = -2451.3333

# Row 22
$ws.Range("H22").Value = 779.1667
$ws.Range("I22").Value = 665
$ws.Range("J22").Value = 893.3333
$ws.Range("K22").Value = 665
$ws.Range("L22").Value = 893.3333
$ws.Range("M22").Value = -370
$ws.Range("N22").Value = -1483.3333

# Row 27
$ws.Range("H27").Value = 779.1667
$ws.Range("I27").Value = 665
$ws.Range("J27").Value = 893.3333
$ws.Range("K27").Value = 665
$ws.Range("L27").Value = 893.3333
$ws.Range("M27").Value = -558
$ws.Range("N27").Value = -1107.3333

# Row 68
$ws.Range("H68").Value = 5399.0527
$ws.Range("I68").Value = 7374.9
$ws.Range("J68").Value = 3203.6667
$ws.Range("K68").Value = 7374.9
$ws.Range("L68").Value = 3203.6667
$ws.Range("M68").Value = -6625.9
$ws.Range("N68").Value = -4701.6667

# Row 71
$ws.Range("H71").Value = 5399.0527
$ws.Range("I71").Value = 7374.9
$ws.Range("J71").Value = 3203.6667
$ws.Range("K71").Value = 36874.5
$ws.Range("L71").Value = 16018.3335
$ws.Range("M71").Value = -33130.5
$ws.Range("N71").Value = -23506.3335

# Row 126
$ws.Range("H126").Value = 3353.75
$ws.Range("I126").Value = 2563.3333
$ws.Range("K126").Value = 7689.999899999999
$ws.Range("M126").Value = -5219.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5308.5
$ws.Range("I62").Value = 5920.4
$ws.Range("J62").Value = 4871.4287
$ws.Range("K62").Value = 5920.4
$ws.Range("L62").Value = 4871.4287
$ws.Range("M62").Value = -5296.4
$ws.Range("N62").Value = -6119.4287

# Row 65
$ws.Range("H65").Value = 5308.5
$ws.Range("I65").Value = 5920.4
$ws.Range("J65").Value = 4871.4287
$ws.Range("K65").Value = 29602
$ws.Range("L65").Value = 24357.1435
$ws.Range("M65").Value = -26482
$ws.Range("N65").Value = -30597.1435
